$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove comments (also drops legacyDrawing relationship) ---
while ($ws.Comments.Count -gt 0) {
    $ws.Comments.Item(1).Delete()
}

# --- Remove hyperlinks on email cells ---
$ws.Range("C2").Hyperlinks.Delete()

# --- Remove the two sample data rows, keep header row only ---
$ws.Rows("2:3").Delete()

# --- Prepare new header cells G1:H1 with same formatting as existing header cells ---
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update header text ---
$ws.Range("E1").Value = "Group"
$ws.Range("F1").Value = "Tags"
$ws.Range("G1").Value = "External Reference ID"
$ws.Range("H1").Value = "External Reference System"

# --- Apply updated header styling (font + fill) across the whole header row ---
$hdr = $ws.Range("A1:H1")
$hdr.Font.ColorIndex = 1
$hdr.Font.Size = 10
$hdr.Font.Name = "Calibri"
$hdr.Interior.PatternColor = 16777215

# --- Column width for Group column ---
$ws.Columns("E").ColumnWidth = 19

# --- Data validation drop-down list for the Group column ---
$ws.Range("E2:E421").Validation.Add(3, 1, 1, '"Group A,Group B,Group C,Group D,Contractual Staff, Others"')
$ws.Range("E2:E421").Validation.InCellDropdown = $true
$ws.Range("E2:E421").Validation.IgnoreBlank = $true

# --- Restore selection similar to the saved workbook ---
$ws.Range("E409").Select()
